# Updates scraped crypto market data (price + 1h volume %) for Sat Nov  2 02:53:23 UTC 2024.
# Source sheet stores every value as literal text (coinranking.com scrape), so numeric-looking
# prices (e.g. "575.17", "1.00") are written with NumberFormat "@" first to stop Excel's COM
# layer from auto-converting them to real numbers (which would silently drop trailing zeros /
# reformat things like "0.0000178"). Percent-change cells keep their original padding spaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.776.09'
$ws.Range('E2').Value = '  +0.70%  '

# Row 3
$ws.Range('D3').Value = '2.515.37'
$ws.Range('E3').Value = '  +0.57%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'  # keep text: '575.17' would otherwise become a Number
$ws.Range('D5').Value = '575.17'
$ws.Range('E5').Value = '  -0.31%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'  # keep text: '167.90' would otherwise become a Number
$ws.Range('D6').Value = '167.90'
$ws.Range('E6').Value = '  +0.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('E8').Value = '  +1.24%  '

# Row 9
$ws.Range('D9').Value = '2.515.52'
$ws.Range('E9').Value = '  +0.67%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'  # keep text: '0.163' would otherwise become a Number
$ws.Range('D10').Value = '0.163'
$ws.Range('E10').Value = '  +1.82%  '

# Row 11
$ws.Range('E11').Value = '  -0.16%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'  # keep text: '0.361' would otherwise become a Number
$ws.Range('D12').Value = '0.361'
$ws.Range('E12').Value = '  +6.78%  '

# Row 13
$ws.Range('E13').Value = '  +1.92%  '

# Row 14
$ws.Range('D14').Value = '2.980.73'
$ws.Range('E14').Value = '  +0.80%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'  # keep text: '0.0000178' would otherwise become a Number
$ws.Range('D15').Value = '0.0000178'
$ws.Range('E15').Value = '  +1.26%  '

# Row 16
$ws.Range('D16').Value = '69.700.48'
$ws.Range('E16').Value = '  +0.67%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'  # keep text: '24.94' would otherwise become a Number
$ws.Range('D17').Value = '24.94'
$ws.Range('E17').Value = '  +0.33%  '

# Row 18
$ws.Range('D18').Value = '2.523.57'
$ws.Range('E18').Value = '  +0.30%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'  # keep text: '11.37' would otherwise become a Number
$ws.Range('D19').Value = '11.37'
$ws.Range('E19').Value = '  +0.18%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'  # keep text: '7.67' would otherwise become a Number
$ws.Range('D20').Value = '7.67'
$ws.Range('E20').Value = '  -0.88%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'  # keep text: '351.90' would otherwise become a Number
$ws.Range('D21').Value = '351.90'
$ws.Range('E21').Value = '  +0.27%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'  # keep text: '3.93' would otherwise become a Number
$ws.Range('D22').Value = '3.93'
$ws.Range('E22').Value = '  +0.02%  '

# Row 23
$ws.Range('E23').Value = '  +0.73%  '

# Row 24
$ws.Range('E24').Value = '  -0.09%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'  # keep text: '70.94' would otherwise become a Number
$ws.Range('D25').Value = '70.94'
$ws.Range('E25').Value = '  +2.98%  '

# Row 26
$ws.Range('E26').Value = '  -0.52%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'  # keep text: '8.89' would otherwise become a Number
$ws.Range('D27').Value = '8.89'
$ws.Range('E27').Value = '  -1.80%  '

# Row 28
$ws.Range('D28').Value = '2.640.58'
$ws.Range('E28').Value = '  +0.39%  '

# Row 29
$ws.Range('E29').Value = '  -0.29%  '

# Row 30
$ws.Range('E30').Value = '  -0.66%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'  # keep text: '7.91' would otherwise become a Number
$ws.Range('D31').Value = '7.91'
$ws.Range('E31').Value = '  +1.05%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'  # keep text: '463.71' would otherwise become a Number
$ws.Range('D32').Value = '463.71'
$ws.Range('E32').Value = '  -4.34%  '

# Row 33
$ws.Range('E33').Value = '  -3.16%  '

# Row 34
$ws.Range('E34').Value = '  +0.06%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'  # keep text: '1.00' would otherwise become a Number
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.21%  '

# Row 36
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'  # keep text: '159.37' would otherwise become a Number
$ws.Range('D36').Value = '159.37'
$ws.Range('E36').Value = '  +2.73%  '

# Row 37
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'  # keep text: '0.117' would otherwise become a Number
$ws.Range('D37').Value = '0.117'
$ws.Range('E37').Value = '  +2.07%  '

# Row 38
$ws.Range('E38').Value = '  +1.19%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'  # keep text: '18.58' would otherwise become a Number
$ws.Range('D39').Value = '18.58'
$ws.Range('E39').Value = '  +0.34%  '

# Row 40
$ws.Range('E40').Value = '  +0.00%  '

# Row 41
$ws.Range('E41').Value = '  +0.71%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'  # keep text: '4.71' would otherwise become a Number
$ws.Range('D42').Value = '4.71'
$ws.Range('E42').Value = '  -0.65%  '

# Row 43
$ws.Range('E43').Value = '  -0.89%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'  # keep text: '38.30' would otherwise become a Number
$ws.Range('D44').Value = '38.30'
$ws.Range('E44').Value = '  +0.12%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'  # keep text: '1.11' would otherwise become a Number
$ws.Range('D45').Value = '1.11'
$ws.Range('E45').Value = '  -5.26%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'  # keep text: '2.23' would otherwise become a Number
$ws.Range('D46').Value = '2.23'
$ws.Range('E46').Value = '  -5.16%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'  # keep text: '143.23' would otherwise become a Number
$ws.Range('D47').Value = '143.23'
$ws.Range('E47').Value = '  +0.56%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'  # keep text: '3.50' would otherwise become a Number
$ws.Range('D48').Value = '3.50'
$ws.Range('E48').Value = '  -0.54%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'  # keep text: '0.524' would otherwise become a Number
$ws.Range('D49').Value = '0.524'
$ws.Range('E49').Value = '  -0.28%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'  # keep text: '0.0737' would otherwise become a Number
$ws.Range('D50').Value = '0.0737'
$ws.Range('E50').Value = '  +1.43%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'  # keep text: '5.82' would otherwise become a Number
$ws.Range('D51').Value = '5.82'
$ws.Range('E51').Value = '  +4.21%  '
